$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update J3:J5 so they mirror the values already held in I3:I5
$ws.Range("J3").Value = $ws.Range("I3").Text
$ws.Range("J4").Value = $ws.Range("I4").Text
$ws.Range("J5").Value = $ws.Range("I5").Text

# Update the selection to match the new active selection (J3:J5, active cell J3)
$ws.Range("J3:J5").Select()
